# Update the Jogos_da_Semana_FlashScore workbook:
#  - tweak a handful of odds on the existing Banfield-Tigre (row 2) and
#    Platense-Godoy Cruz (row 4) matches
#  - insert a brand-new match (Botafogo SP vs Avai, Brazil - Serie B) as
#    the new row 5, pushing the two Uruguay matches down one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Banfield vs Tigre): a few odds changed ---
$ws.Range("J2").Value = 3.4
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("BD2").Value = 126

# --- Row 4 (Platense vs Godoy Cruz): several odds changed ---
$ws.Range("H4").Value = 3.1
$ws.Range("J4").Value = 2.75
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("S4").Value = 1.73
$ws.Range("T4").Value = 2.08
$ws.Range("AE4").Value = 26
$ws.Range("AN4").Value = 3.5
$ws.Range("AU4").Value = 11
$ws.Range("AW4").Value = 6.5

# --- Insert a brand new row at position 5 for the new match ---
$ws.Rows.Item(5).Insert()

$newRow = @{
    "A" = "nyeDPXV6"
    "B" = "18/11/2024"
    "C" = "21:00"
    "D" = "BRAZIL - SERIE B"
    "E" = "Botafogo SP"
    "F" = "Avai"
    "G" = 2.5
    "H" = 2.9
    "I" = 3
    "J" = 3.5
    "K" = 1.83
    "L" = 4
    "M" = 1.13
    "N" = 6
    "O" = 1.62
    "P" = 2.2
    "Q" = 3.1
    "R" = 1.36
    "S" = 1.62
    "T" = 2.2
    "U" = 2.25
    "V" = 1.57
    "W" = 6
    "X" = 10
    "Y" = 11
    "Z" = 26
    "AA" = 26
    "AB" = 41
    "AC" = 6
    "AD" = 6
    "AE" = 21
    "AF" = 81
    "AG" = 6.5
    "AH" = 13
    "AI" = 13
    "AJ" = 34
    "AK" = 34
    "AL" = 51
    "AM" = 201
    "AN" = 4.33
    "AO" = 17
    "AP" = 34
    "AQ" = 51
    "AR" = 101
    "AS" = 351
    "AT" = 2.2
    "AU" = 9.5
    "AV" = 81
    "AW" = 4.75
    "AX" = 19
    "AY" = 34
    "AZ" = 67
    "BA" = 126
    "BB" = 500
    "BC" = 81
    "BD" = 81
}

foreach ($col in $newRow.Keys) {
    $ws.Range("$col" + "5").Value = $newRow[$col]
}

# --- Row 7 (the old CA Cerro vs Boston River row, now shifted down) also
#     got a handful of odds re-priced in this update ---
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2.2
$ws.Range("R7").Value = 1.65
